$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (was "BrassA-HW20.xpc")
$ws.Name = "BrassA"

# Tiny precision corrections in row 13 (re-computed with the new
# Gaussian-quadrature routine)
$ws.Range("D13").Value = 0.9980247514947538
$ws.Range("H13").Value = 0.9980247514947538
$ws.Range("N13").Value = 0.9953297297321173
$ws.Range("O13").Value = 0.9962326804467529

# Append a new data row (14 / HexGrid-60degTilt5degRes) exported from the
# Gaussian Quadrature Scheme, mirroring the formatting of the row above it
$ws.Range("A15:P15").Copy($ws.Range("A16"))

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.187930224880939
$ws.Range("D16").Value = 0.9440458729845717
$ws.Range("E16").Value = 0.9527852964205011
$ws.Range("F16").Value = 0.9626495799918713
$ws.Range("G16").Value = 1.187930224880939
$ws.Range("H16").Value = 0.9440458729845717
$ws.Range("I16").Value = 1.022882622533515
$ws.Range("J16").Value = 0.9260464903722608
$ws.Range("K16").Value = 1.045717215343752
$ws.Range("L16").Value = 0.9352641577527503
$ws.Range("M16").Value = 1.187930224880939
$ws.Range("N16").Value = 0.9484155847025364
$ws.Range("O16").Value = 1.011852743569471
$ws.Range("P16").Value = 0.9971651825350201
